$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update genet names for merging with 2020 data (append "RR" suffix)
$ws.Range("D55").Value = "Acer121RR"
$ws.Range("D59").Value = "Acer120RR"
$ws.Range("D67").Value = "Acer123RR"

# Fix genet names: update cell values in column D (genet).
# ML71 should have been ML77 (fix per field notes)
$ws.Range("D35").Value = "ML77"

# Update view state: selection position when saved
$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F35").Select()

$wb.Save()
